$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4666194915771484
$ws.Range("E2").Value = 207.7360586558789
$ws.Range("F2").Value = 0.005953584294917272
$ws.Range("G2").Value = 0.005368189131073895
$ws.Range("H2").Value = 0.005098339670967602
$ws.Range("I2").Value = 0.004941235681427699
$ws.Range("J2").Value = 0.004867755711073091
$ws.Range("K2").Value = 0.004489686657707564
$ws.Range("L2").Value = 0.004489686657707564
$ws.Range("M2").Value = 0.004402221308805502
$ws.Range("N2").Value = 0.004216176235565184
$ws.Range("O2").Value = 0.004216176235565184
$ws.Range("P2").Value = 0.004202602549725851
$ws.Range("Q2").Value = 0.004202602549725851
$ws.Range("R2").Value = 0.004202602549725851
$ws.Range("S2").Value = 0.004151112666490322
$ws.Range("T2").Value = 0.004145101186604272
$ws.Range("U2").Value = 0.004137163923682215
$ws.Range("V2").Value = 0.004104204417666121
$ws.Range("W2").Value = 0.004083757561946301
$ws.Range("X2").Value = 0.004080544764813053
$ws.Range("Y2").Value = 0.004049435841245201

$ws.Range("C3").Value = 0.4045958518981934
$ws.Range("E3").Value = 215.9554006296148
$ws.Range("F3").Value = 0.006012554821459275
$ws.Range("G3").Value = 0.005295036367631936
$ws.Range("H3").Value = 0.004986953847832086
$ws.Range("I3").Value = 0.004895761688391418
$ws.Range("J3").Value = 0.004873400258503559
$ws.Range("K3").Value = 0.004576401779807722
$ws.Range("L3").Value = 0.004576401779807722
$ws.Range("M3").Value = 0.004576401779807722
$ws.Range("N3").Value = 0.00449700150891417
$ws.Range("O3").Value = 0.004388778306368401
$ws.Range("P3").Value = 0.004353924678899242
$ws.Range("Q3").Value = 0.004353924678899242
$ws.Range("R3").Value = 0.004332498008071731
$ws.Range("S3").Value = 0.004318603346537006
$ws.Range("T3").Value = 0.004260733166965895
$ws.Range("U3").Value = 0.004246970450866192
$ws.Range("V3").Value = 0.0042404825388545
$ws.Range("W3").Value = 0.004227657068647609
$ws.Range("X3").Value = 0.004221391626172963
$ws.Range("Y3").Value = 0.004209656932351165

$ws.Range("C4").Value = 0.4070723056793213
$ws.Range("E4").Value = 216.3366982265979
$ws.Range("F4").Value = 0.006067332779993517
$ws.Range("G4").Value = 0.005393267121205752
$ws.Range("H4").Value = 0.004943868988933328
$ws.Range("I4").Value = 0.004901353057935481
$ws.Range("J4").Value = 0.004783275302654206
$ws.Range("K4").Value = 0.004725725655377907
$ws.Range("L4").Value = 0.004683829221917478
$ws.Range("M4").Value = 0.004509412023237326
$ws.Range("N4").Value = 0.004509412023237326
$ws.Range("O4").Value = 0.00449572647034645
$ws.Range("P4").Value = 0.004349988861229302
$ws.Range("Q4").Value = 0.004349988861229302
$ws.Range("R4").Value = 0.004349988861229302
$ws.Range("S4").Value = 0.004318443919448759
$ws.Range("T4").Value = 0.004292050975768217
$ws.Range("U4").Value = 0.004247510690617281
$ws.Range("V4").Value = 0.004247510690617281
$ws.Range("W4").Value = 0.004235761174453334
$ws.Range("X4").Value = 0.004226935847431586
$ws.Range("Y4").Value = 0.004217089634046742

$ws.Range("C5").Value = 0.3678033351898193
$ws.Range("E5").Value = 210.4043682009433
$ws.Range("F5").Value = 0.005989927846952334
$ws.Range("G5").Value = 0.005359822080949472
$ws.Range("H5").Value = 0.004992930339307043
$ws.Range("I5").Value = 0.004626250346772948
$ws.Range("J5").Value = 0.004626250346772948
$ws.Range("K5").Value = 0.004626250346772948
$ws.Range("L5").Value = 0.004535987767166975
$ws.Range("M5").Value = 0.004472027091651575
$ws.Range("N5").Value = 0.00425663050617355
$ws.Range("O5").Value = 0.00425663050617355
$ws.Range("P5").Value = 0.00425663050617355
$ws.Range("Q5").Value = 0.00425663050617355
$ws.Range("R5").Value = 0.00425663050617355
$ws.Range("S5").Value = 0.004238646207718029
$ws.Range("T5").Value = 0.00422082802008896
$ws.Range("U5").Value = 0.00422082802008896
$ws.Range("V5").Value = 0.004193705817433935
$ws.Range("W5").Value = 0.004138815701934746
$ws.Range("X5").Value = 0.004101449672533007
$ws.Range("Y5").Value = 0.004101449672533007

$ws.Range("C6").Value = 0.390594482421875
$ws.Range("E6").Value = 213.2684369622966
$ws.Range("F6").Value = 0.005940771258625728
$ws.Range("G6").Value = 0.0051573911929168
$ws.Range("H6").Value = 0.004853764315462542
$ws.Range("I6").Value = 0.004853764315462542
$ws.Range("J6").Value = 0.00460177096071283
$ws.Range("K6").Value = 0.00460177096071283
$ws.Range("L6").Value = 0.004492357454156328
$ws.Range("M6").Value = 0.004386181290035708
$ws.Range("N6").Value = 0.004386181290035708
$ws.Range("O6").Value = 0.004352352239860082
$ws.Range("P6").Value = 0.00432182339647724
$ws.Range("Q6").Value = 0.00432182339647724
$ws.Range("R6").Value = 0.004284707444585754
$ws.Range("S6").Value = 0.004278750701051873
$ws.Range("T6").Value = 0.00421469357416699
$ws.Range("U6").Value = 0.004198743302276746
$ws.Range("V6").Value = 0.004187467340536996
$ws.Range("W6").Value = 0.004180964714938121
$ws.Range("X6").Value = 0.004177106434668267
$ws.Range("Y6").Value = 0.004157279472949251

$ws.Range("C7").Value = 0.4305441379547119
$ws.Range("E7").Value = 212.4196194185279
$ws.Range("F7").Value = 0.005961684840401669
$ws.Range("G7").Value = 0.005240373657493
$ws.Range("H7").Value = 0.004968727956311801
$ws.Range("I7").Value = 0.004915811333706543
$ws.Range("J7").Value = 0.004915811333706543
$ws.Range("K7").Value = 0.004852708101142919
$ws.Range("L7").Value = 0.004687305311851096
$ws.Range("M7").Value = 0.004687305311851096
$ws.Range("N7").Value = 0.004609103805919198
$ws.Range("O7").Value = 0.004417282740778778
$ws.Range("P7").Value = 0.004417282740778778
$ws.Range("Q7").Value = 0.004379634664962038
$ws.Range("R7").Value = 0.004363428300196522
$ws.Range("S7").Value = 0.004343192812140476
$ws.Range("T7").Value = 0.004296120684679014
$ws.Range("U7").Value = 0.004222417431783735
$ws.Range("V7").Value = 0.00420004693556851
$ws.Range("W7").Value = 0.004169458281882333
$ws.Range("X7").Value = 0.004155395403456503
$ws.Range("Y7").Value = 0.004140733321998593

$ws.Range("C8").Value = 0.4419827461242676
$ws.Range("E8").Value = 208.8545051639248
$ws.Range("F8").Value = 0.005879060393559895
$ws.Range("G8").Value = 0.005123601416366049
$ws.Range("H8").Value = 0.005085588020634996
$ws.Range("I8").Value = 0.004824729987737871
$ws.Range("J8").Value = 0.004765616560080213
$ws.Range("K8").Value = 0.004624287067181035
$ws.Range("L8").Value = 0.004527822059131601
$ws.Range("M8").Value = 0.004484219655513898
$ws.Range("N8").Value = 0.004465508577135618
$ws.Range("O8").Value = 0.004333292069539055
$ws.Range("P8").Value = 0.00429474258212233
$ws.Range("Q8").Value = 0.004190986837743765
$ws.Range("R8").Value = 0.004190986837743765
$ws.Range("S8").Value = 0.004190986837743765
$ws.Range("T8").Value = 0.004177371684374527
$ws.Range("U8").Value = 0.004136219046032285
$ws.Range("V8").Value = 0.004131651043263351
$ws.Range("W8").Value = 0.004106699945534191
$ws.Range("X8").Value = 0.004106699945534191
$ws.Range("Y8").Value = 0.004071237917425434

$ws.Range("C9").Value = 0.3808796405792236
$ws.Range("E9").Value = 212.2347159621313
$ws.Range("F9").Value = 0.005912705158227065
$ws.Range("G9").Value = 0.005022365187489082
$ws.Range("H9").Value = 0.004968408541628756
$ws.Range("I9").Value = 0.004968408541628756
$ws.Range("J9").Value = 0.004508968004553897
$ws.Range("K9").Value = 0.004508968004553897
$ws.Range("L9").Value = 0.004478828371294852
$ws.Range("M9").Value = 0.004449892808438464
$ws.Range("N9").Value = 0.004449892808438464
$ws.Range("O9").Value = 0.00433980375303997
$ws.Range("P9").Value = 0.00433980375303997
$ws.Range("Q9").Value = 0.004308410226834323
$ws.Range("R9").Value = 0.004280522964004082
$ws.Range("S9").Value = 0.004280522964004082
$ws.Range("T9").Value = 0.004277811651051402
$ws.Range("U9").Value = 0.004219692011902863
$ws.Range("V9").Value = 0.004193212839220055
$ws.Range("W9").Value = 0.004190584126509884
$ws.Range("X9").Value = 0.004149368291141372
$ws.Range("Y9").Value = 0.004137128966123417

$ws.Range("C10").Value = 0.3906245231628418
$ws.Range("E10").Value = 220.5958513863643
$ws.Range("F10").Value = 0.006030707666437985
$ws.Range("G10").Value = 0.005399849135774593
$ws.Range("H10").Value = 0.004974936891654297
$ws.Range("I10").Value = 0.004834682809387745
$ws.Range("J10").Value = 0.004776794918418114
$ws.Range("K10").Value = 0.004752282765899688
$ws.Range("L10").Value = 0.004752282765899688
$ws.Range("M10").Value = 0.004700939901956223
$ws.Range("N10").Value = 0.004700939901956223
$ws.Range("O10").Value = 0.004636767666316623
$ws.Range("P10").Value = 0.004616372323600397
$ws.Range("Q10").Value = 0.004604363258842521
$ws.Range("R10").Value = 0.004525373939545422
$ws.Range("S10").Value = 0.004466050548733186
$ws.Range("T10").Value = 0.004438455446793917
$ws.Range("U10").Value = 0.00435752736712277
$ws.Range("V10").Value = 0.004324757431058366
$ws.Range("W10").Value = 0.004320299179355705
$ws.Range("X10").Value = 0.004300114062112364
$ws.Range("Y10").Value = 0.004300114062112364

$ws.Range("C11").Value = 0.3906254768371582
$ws.Range("E11").Value = 213.1421084955982
$ws.Range("F11").Value = 0.006051709744267044
$ws.Range("G11").Value = 0.005371799740082598
$ws.Range("H11").Value = 0.005077553426418004
$ws.Range("I11").Value = 0.004946435328762857
$ws.Range("J11").Value = 0.004773064965607439
$ws.Range("K11").Value = 0.004773064965607439
$ws.Range("L11").Value = 0.004579564551253032
$ws.Range("M11").Value = 0.00452281300209219
$ws.Range("N11").Value = 0.00452281300209219
$ws.Range("O11").Value = 0.004269404317556658
$ws.Range("P11").Value = 0.004269404317556658
$ws.Range("Q11").Value = 0.004269404317556658
$ws.Range("R11").Value = 0.004269404317556658
$ws.Range("S11").Value = 0.004212804742222211
$ws.Range("T11").Value = 0.004212804742222211
$ws.Range("U11").Value = 0.004212804742222211
$ws.Range("V11").Value = 0.004212359845887103
$ws.Range("W11").Value = 0.004170835118506107
$ws.Range("X11").Value = 0.004170835118506107
$ws.Range("Y11").Value = 0.004154816929738756
